$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A106").Value = "$ 27.617 CLP 31-12-20"
$ws.Range("A107").Value = "$ 27.617 CLP 31-12-20"
$ws.Range("A108").Value = "$ 27.613 CLP 04-01-21"
$ws.Range("A109").Value = "$ 27.613 CLP 04-01-21"
$ws.Range("A110").Value = "$ 34.589 CLP 04-01-21"
$ws.Range("A111").Value = "$ 34.589 CLP 04-01-21"
$ws.Range("A112").Value = "0.95 UF 13-01-21"
$ws.Range("A113").Value = "0.95 UF 13-01-21"
